# ---- Edit: add RootNamespace + AssemblyName mapping sheets, ----
# ---- update Maui.Controls B47, and adjust selections       ----
$wb = $excel.ActiveWorkbook

# 1) Update existing Maui.Controls sheet: Xamarin.Flex -> Microsoft.Maui.Controls.Flex
$wsControls = $wb.Worksheets.Item("Maui.Controls")
$wsControls.Range("B47").Value = "Microsoft.Maui.Controls.Flex"
$wsControls.Range("B48").Select()

# 2) Files sheet selection moves (no longer the tab-selected sheet)
$wsFiles = $wb.Worksheets.Item("Files")
$wsFiles.Range("A58").Select()

# 3) Add the RootNamespace sheet at the end
$count = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($count)
$wsRoot = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$wsRoot.Name = "RootNamespace"
$wsRoot.Range("A1:B1").ColumnWidth = 51.8

$rootData = @(
  @("Xamarin.Forms.Platform.Android", "Microsoft.Maui.Controls.Compatibility.Android"),
  @("Xamarin.Forms.Platform.Android.AppLinks", "Microsoft.Maui.Controls.Compatibility.Android.AppLinks"),
  @("Xamarin.Forms.Platform.GTK", "Microsoft.Maui.Controls.Compatibility.GTK"),
  @("Xamarin.Forms.Platform.iOS", "Microsoft.Maui.Controls.Compatibility.iOS"),
  @("Xamarin.Forms.Platform.macOS", "Microsoft.Maui.Controls.Compatibility.macOS"),
  @("Xamarin.Forms.Platform.UAP", "Microsoft.Maui.Controls.Compatibility.UAP"),
  @("Xamarin.Forms.Platform.Android.UnitTests", "Microsoft.Maui.Controls.Compatibility.Android.UnitTests"),
  @("Xamarin.Forms.Platform.iOS.UnitTests", "Microsoft.Maui.Controls.Compatibility.iOS.UnitTests"),
  @("Xamarin.Forms.Platform.UAP.UnitTests", "Microsoft.Maui.Controls.Compatibility.UAP.UnitTests"),
  @("Xamarin.Forms.Maps.Android", "Microsoft.Maui.Controls.Compatibility.Android"),
  @("Xamarin.Forms.Maps.GTK", "Microsoft.Maui.Controls.GTK"),
  @("Xamarin.Forms.Maps.iOS", "Microsoft.Maui.Controls.iOS"),
  @("Xamarin.Forms.Maps.MacOS", "Microsoft.Maui.Controls.MacOS"),
  @("Xamarin.Forms.Maps.UWP", "Microsoft.Maui.Controls.UWP"),
  @("Xamarin.Forms.Material.Android", "Microsoft.Maui.Controls.Compatibility.Material.Android"),
  @("Xamarin.Forms.Material.iOS", "Microsoft.Maui.Controls.Compatibility.Material.iOS"),
  @("Xamarin.Forms.Material.Tizen", "Microsoft.Maui.Controls.Compatibility.Material.Tizen"),
  @("Xamarin.Platform.Handlers.DeviceTests", "Microsoft.Maui.Handlers.DeviceTests"),
  @("Xamarin.Platform", "Microsoft.Maui"),
  @("Xamarin.Essentials", "Microsoft.Maui.Essentials"),
  @("Xamarin.Forms.Xaml.Design", "Microsoft.Maui.Controls.Xaml.Design"),
  @("Xamarin.Forms.Xaml.UnitTests", "Microsoft.Maui.Controls.Xaml.UnitTests"),
  @("Xamarin.Forms.DualScreen", "Microsoft.Maui.Controls.DualScreen"),
  @("Xamarin.Forms.Core", "Microsoft.Maui.Controls.Core")
)

$r = 2
foreach ($pair in $rootData) {
    $wsRoot.Cells.Item($r, 1).Value = $pair[0]
    $wsRoot.Cells.Item($r, 2).Value = $pair[1]
    $r = $r + 1
}
$wsRoot.Range("B26").Select()

# 4) Add the AssemblyName sheet at the end (this becomes the active/selected sheet)
$count = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($count)
$wsAsm = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$wsAsm.Name = "AssemblyName"
$wsAsm.Range("A1:B1").ColumnWidth = 56.2

$asmData = @(
  @("Xamarin.Forms.Platform.Android", "Microsoft.Maui.Controls.Compatibility.Android"),
  @("Xamarin.Forms.Platform.Android.AppLinks", "Microsoft.Maui.Controls.Compatibility.Android.AppLinks"),
  @("Xamarin.Forms.Platform.GTK", "Microsoft.Maui.Controls.Compatibility.GTK"),
  @("Xamarin.Forms.Platform.iOS", "Microsoft.Maui.Controls.Compatibility.iOS"),
  @("Xamarin.Forms.Platform.macOS", "Microsoft.Maui.Controls.Compatibility.macOS"),
  @("Xamarin.Forms.Platform.Android.UnitTests", "Microsoft.Maui.Controls.Compatibility.Android.UnitTests"),
  @("Xamarin.Forms.Platform.iOS.UnitTests", "Microsoft.Maui.Controls.Compatibility.iOS.UnitTests"),
  @("Xamarin.Forms.Platform.UAP.UnitTests", "Microsoft.Maui.Controls.Compatibility.UAP.UnitTests"),
  @("Xamarin.Forms.Maps.Android", "Microsoft.Maui.Controls.Maps.Android"),
  @("Xamarin.Forms.Maps.GTK", "Microsoft.Maui.Controls.Maps.GTK"),
  @("Xamarin.Forms.Maps.iOS", "Microsoft.Maui.Controls.Maps.iOS"),
  @("Xamarin.Forms.Maps.macOS", "Microsoft.Maui.Controls.Maps.macOS"),
  @("Xamarin.Forms.Maps.UWP", "Microsoft.Maui.Controls.Maps.UWP"),
  @("Xamarin.Forms.Material", "Microsoft.Maui.Controls.Compatibility.Material"),
  @("Xamarin.Forms.Build.Tasks", "Microsoft.Maui.Controls.Build.Tasks"),
  @("Xamarin.Forms.Core.Design", "Microsoft.Maui.Controls.Core.Design"),
  @("Xamarin.Forms.DualScreen", "Microsoft.Maui.Controls.DualScreen"),
  @("Xamarin.Forms.Xaml.UnitTests", "Microsoft.Maui.Controls.Xaml.UnitTests"),
  @("Xamarin.Forms.Xaml.Design", "Microsoft.Maui.Controls.Xaml.Design"),
  @("Xamarin.Essentials", "Microsoft.Maui.Essentials"),
  @("Xamarin.Forms", "Microsoft.Maui.Controls"),
  @("XamarinEssentialsDeviceTestsiOS", "EssentialsDeviceTestsiOS"),
  @("XamarinEssentialsDeviceTestsAndroid", "EssentialsDeviceTestsAndroid"),
  @("XamarinEssentialsDeviceTestsShared", "EssentialsDeviceTestsShared"),
  @("XamarinEssentialsDeviceTestsUWP", "EssentialsDeviceTestsUWP"),
  @("XamarinEssentialsTests", "EssentialsTests"),
  @("Xamarin.Platform.Handlers.DeviceTests", "Microsoft.Maui.Handlers.DeviceTests")
)

$r = 2
foreach ($pair in $asmData) {
    $wsAsm.Cells.Item($r, 1).Value = $pair[0]
    $wsAsm.Cells.Item($r, 2).Value = $pair[1]
    $r = $r + 1
}
$wsAsm.Range("B15").Select()
$wsAsm.Activate()

Write-Output "edit complete"
